$d = $word.ActiveDocument

# --- Paragraph "Clases y objetos:" ---
# The trailing two runs ("... creación de clases y objetos, " and
# "encapsulación, herencia y polimorfismo.") lose their yellow highlight and
# get merged into a single run once the highlighting difference disappears.
$combinedText = " Conceptos de programación orientada a objetos (POO), creación de clases y objetos, encapsulación, herencia y polimorfismo."
$full1 = "Clases y objetos:" + $combinedText

$r = $d.Content
$r.Find.Execute($full1, $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.HighlightColorIndex = 0          # wdNoHighlight -> strips run-level highlight
$r.Collapse(0)                      # collapse to the paragraph mark
$r.Font.HighlightColorIndex = 0     # strips the paragraph-mark's own highlight

# Merge the two now-identically-formatted runs into a single run via a
# self find-and-replace over just their shared text span.
$d.Content.Find.Execute($combinedText, $false, $false, $false, $false, $false, $true, 1, $false, $combinedText, 2) | Out-Null

# --- Paragraph "Manejo de excepciones en POO:" ---
$full3 = "Manejo de excepciones en POO: Uso de try, except, finally y la creación de excepciones personalizadas."
$r3 = $d.Content
$r3.Find.Execute($full3, $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r3.HighlightColorIndex = 0
$r3.Collapse(0)
$r3.Font.HighlightColorIndex = 0

# --- Paragraph "Patrones de diseño básicos en Python:" ---
$full4 = "Patrones de diseño básicos en Python: Singleton, Factory, y Strategy aplicados a la POO."
$r4 = $d.Content
$r4.Find.Execute($full4, $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r4.HighlightColorIndex = 0
$r4.Collapse(0)
$r4.Font.HighlightColorIndex = 0
